# Rimappa le righe 2..27 in base al nuovo ordinamento dei gruppi "macchina"
# (le funzioni di ricerca locale ora usano insert invece di move).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 27
$lastCol = 19   # colonne A..S

# Mappa: riga di origine (prima) -> riga di destinazione (dopo)
$rowMap = @{
    2  = 5
    3  = 6
    4  = 7
    5  = 8
    6  = 9
    7  = 10
    8  = 11
    9  = 12
    10 = 13
    11 = 15
    12 = 16
    13 = 17
    14 = 27
    15 = 14
    16 = 2
    17 = 3
    18 = 4
    19 = 18
    20 = 19
    21 = 20
    22 = 21
    23 = 22
    24 = 23
    25 = 24
    26 = 25
    27 = 26
}

# 1) Istantanea di tutti i valori correnti (prima che vengano sovrascritti)
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowVals = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $rowVals += $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowVals
}

# 2) Scrive ogni riga nella nuova posizione in base alla mappa
foreach ($srcRow in $rowMap.Keys) {
    $dstRow = $rowMap[$srcRow]
    $vals = $snapshot[$srcRow]
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($dstRow, $c).Value2 = $vals[$c - 1]
    }
}
